$d = $word.ActiveDocument

# --- Edit 1: strike through "reinterpret_" in the trunc_/convert_/... list ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text
    if ($full -like "*reinterpret_ functions*") {
        $pStart = $p.Range.Start
        $idx = $full.IndexOf("reinterpret_")
        $start = $pStart + $idx
        $end = $start + "reinterpret_".Length
        $sub = $d.Range($start, $end)
        $sub.Font.StrikeThrough = 1
        break
    }
}

# --- Edit 2: strike through "memory.fill, memory.copy," in the memory export list ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text
    if ($full -like "*memory.fill*memory.copy*memory.init*") {
        $pStart = $p.Range.Start
        $idxStart = $full.IndexOf("memory.fill")
        $idxEnd = $full.IndexOf("memory.init")
        $start = $pStart + $idxStart
        $end = $pStart + $idxEnd - 1
        $sub = $d.Range($start, $end)
        $sub.Font.StrikeThrough = 1
        break
    }
}
